# Commit: Renamed "Shrub" to "CSS"
#
# The Tukey posthoc table has a single comparison row (row 2):
#   group1 / group2 / meandiff / p-adj / lower / upper / reject
# Originally: Grassland vs Shrub, meandiff=-2.3356, CI=[-2.9039, -1.7672]
#
# The "Shrub" vegetation treatment is renamed to "CSS". Because group1/group2
# are sorted alphabetically, the pair becomes CSS vs Grassland, and the
# meandiff/CI bounds are negated (and lower/upper swapped) to reflect the
# swapped subtraction order:
#   New: CSS vs Grassland, meandiff=2.3356, CI=[1.7672, 2.9039]

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# group1 / group2 labels
$ws.Range("A2").Value = "CSS"
$ws.Range("B2").Value = "Grassland"

# meandiff, lower, upper (p-adj and reject stay the same)
$ws.Range("C2").Value = 2.3356
$ws.Range("E2").Value = 1.7672
$ws.Range("F2").Value = 2.9039
